$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C4").Value = -12.01059999999999
$ws.Range("A9").Value = -21.82159999999999
$ws.Range("C9").Value = -12.6888
$ws.Range("D9").Value = -8.678500000000001
$ws.Range("C11").Value = -12.9482
$ws.Range("A18").Value = -22.22600000000001
$ws.Range("A20").Value = -21.08769999999999
$ws.Range("C23").Value = -12.1333
$ws.Range("C24").Value = -13.48669999999999
$ws.Range("C26").Value = -12.56060000000001
$ws.Range("A27").Value = -21.8577
$ws.Range("D27").Value = -7.941100000000002
$ws.Range("D29").Value = -7.001300000000001
$ws.Range("D32").Value = -6.946699999999993
$ws.Range("C34").Value = -12.31080000000001
$ws.Range("A35").Value = -21.38449999999999
$ws.Range("C35").Value = -11.6007
$ws.Range("D37").Value = -7.159100000000002
$ws.Range("D38").Value = -8.199399999999992
$ws.Range("D41").Value = -7.639699999999999
$ws.Range("D45").Value = -7.031999999999997
$ws.Range("C48").Value = -10.9064
$ws.Range("C49").Value = -13.8905
$ws.Range("D51").Value = -8.441099999999999
$ws.Range("C52").Value = -10.92099999999999
$ws.Range("D57").Value = -8.205599999999999
$ws.Range("D64").Value = -7.44969999999999
$ws.Range("C66").Value = -10.9446
$ws.Range("C67").Value = -10.84979999999999
$ws.Range("A69").Value = -21.6602
$ws.Range("A76").Value = -19.28409999999998
$ws.Range("A78").Value = -20.33939999999998
$ws.Range("C78").Value = -12.52740000000001
$ws.Range("C80").Value = -13.0338
$ws.Range("A82").Value = -22.10380000000003
$ws.Range("D82").Value = -8.629900000000008
$ws.Range("A83").Value = -21.51989999999999
$ws.Range("A93").Value = -21.44790000000002
$ws.Range("D93").Value = -6.743599999999995
$ws.Range("C99").Value = -13.0131
$ws.Range("D102").Value = -7.212
$ws.Range("C104").Value = -12.6258
$ws.Range("D105").Value = -7.499199999999999
